$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update prices for rows 7 and 8
$ws.Range("C7").Value = 16377
$ws.Range("C8").Value = 15890

# Update shop links (F9, F8, F7) to the new Nix motherboard URLs, keeping
# the existing hyperlink cell style intact.
$asusUrl = "https://www.nix.ru/autocatalog/motherboards_asustek/ASUS-PRIME-B660-PLUS-D4-RTL-LGA1700-B660-2xPCI-E-Dsub-plus-HDMI-plus-DP-25GbLAN-SATA-ATX-4DDR4_574498.html"
$msiUrl = "https://www.nix.ru/autocatalog/motherboards_msi/MSI-MAG-B550-TOMAHAWK-MAX-WIFI-RTL-AM4-B550-2xPCI-E-HDMI-plus-DP25GbLAN-plus-WiFi-plus-BT-SATA-ATX-4DDR4_679259.html"

$styleF9 = $ws.Range("F9").Style
$ws.Hyperlinks.Add($ws.Range("F9"), $asusUrl, "", "", $asusUrl)
$ws.Range("F9").Style = $styleF9

$styleF8 = $ws.Range("F8").Style
$ws.Hyperlinks.Add($ws.Range("F8"), $asusUrl, "", "", $asusUrl)
$ws.Range("F8").Style = $styleF8

$styleF7 = $ws.Range("F7").Style
$ws.Hyperlinks.Add($ws.Range("F7"), $msiUrl, "", "", $msiUrl)
$ws.Range("F7").Style = $styleF7

# Move the active selection to C8 (matches the saved window state)
$ws.Range("C8").Select() | Out-Null
